$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before the old row 177 (weekly data refresh:
# the previously-last two entries bump down and two fresh weekly rows
# are added at the top of this date block).
$ws.Rows("177:178").Insert()

# New row 177 (Segunda)
$ws.Range("A177").Value = 1
$ws.Range("B177").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C177").Value = "Arica y Parinacota"
$ws.Range("D177").Value = 44448
$ws.Range("E177").Value = 15
$ws.Range("F177").Value = 100112023
$ws.Range("G177").Value = "Brócoli"
$ws.Range("H177").Value = "Sin especificar"
$ws.Range("I177").Value = "Segunda"
$ws.Range("J177").Value = 1000
$ws.Range("K177").Value = 800
$ws.Range("L177").Value = 900
$ws.Range("M177").Value = 850
$ws.Range("N177").Value = "`$/unidad"
$ws.Range("O177").Value = "Región de Arica y Parinacota"
$ws.Range("P177").Value = 850
$ws.Range("Q177").Value = 1
$ws.Range("R177").Value = "Hortaliza"

# New row 178 (Tercera)
$ws.Range("A178").Value = 1
$ws.Range("B178").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C178").Value = "Arica y Parinacota"
$ws.Range("D178").Value = 44448
$ws.Range("E178").Value = 15
$ws.Range("F178").Value = 100112023
$ws.Range("G178").Value = "Brócoli"
$ws.Range("H178").Value = "Sin especificar"
$ws.Range("I178").Value = "Tercera"
$ws.Range("J178").Value = 700
$ws.Range("K178").Value = 600
$ws.Range("L178").Value = 700
$ws.Range("M178").Value = 650
$ws.Range("N178").Value = "`$/unidad"
$ws.Range("O178").Value = "Región de Arica y Parinacota"
$ws.Range("P178").Value = 650
$ws.Range("Q178").Value = 1
$ws.Range("R178").Value = "Hortaliza"
